# Updated symbol list on Mon Jan 23 05:43:58 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# coin rows whose quotes changed. Values are written with a leading
# apostrophe so Excel stores them as literal text (matching the workbook's
# existing inline-string cells) instead of auto-converting "305.55" to a
# number or "1.30%" to a percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.55"
$ws.Range("E2").Value = "'1.30%"
$ws.Range("D3").Value = "'36.13"
$ws.Range("E3").Value = "'-4.15%"
$ws.Range("D4").Value = "'5.072"
$ws.Range("E4").Value = "'1.24%"
$ws.Range("D5").Value = "'0.07863"
$ws.Range("E5").Value = "'0.01%"
$ws.Range("D6").Value = "'2.181"
$ws.Range("E6").Value = "'-0.28%"
$ws.Range("D7").Value = "'7.919"
$ws.Range("E7").Value = "'-1.06%"
$ws.Range("D8").Value = "'0.9188"
$ws.Range("D9").Value = "'0.09738"
$ws.Range("E9").Value = "'5.50%"
$ws.Range("D10").Value = "'0.1861"
$ws.Range("E10").Value = "'-0.55%"
$ws.Range("D11").Value = "'0.08637"
$ws.Range("E11").Value = "'2.02%"
$ws.Range("D12").Value = "'0.03497"
$ws.Range("E12").Value = "'-0.48%"
$ws.Range("D13").Value = "'0.09934"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("D14").Value = "'0.001445"
$ws.Range("E14").Value = "'-1.74%"
$ws.Range("D15").Value = "'0.005632"
$ws.Range("E15").Value = "'0.15%"
$ws.Range("E16").Value = "'-0.49%"
$ws.Range("D17").Value = "'4.094"
$ws.Range("E17").Value = "'2.26%"
$ws.Range("D18").Value = "'2.417"
$ws.Range("E18").Value = "'15.36%"
$ws.Range("D19").Value = "'0.3427"
$ws.Range("E19").Value = "'-1.06%"
$ws.Range("D20").Value = "'0.1357"
$ws.Range("E20").Value = "'3.11%"
$ws.Range("D21").Value = "'4.816"
$ws.Range("E21").Value = "'5.87%"
$ws.Range("D22").Value = "'0.2209"
$ws.Range("E22").Value = "'-1.43%"
$ws.Range("D23").Value = "'0.04557"
$ws.Range("E23").Value = "'-1.86%"
$ws.Range("D24").Value = "'0.005093"
$ws.Range("E24").Value = "'14.46%"
$ws.Range("E25").Value = "'0.50%"
$ws.Range("D26").Value = "'0.0001400"
$ws.Range("E26").Value = "'7.92%"
$ws.Range("D27").Value = "'0.0004757"
$ws.Range("E27").Value = "'0.32%"
$ws.Range("D39").Value = "'0.01826"
$ws.Range("E39").Value = "'4.58%"
$ws.Range("D40").Value = "'0.04764"
$ws.Range("E40").Value = "'1.06%"
$ws.Range("D41").Value = "'0.007733"
$ws.Range("E41").Value = "'-2.14%"
$ws.Range("D42").Value = "'0.1396"
$ws.Range("E42").Value = "'0.39%"
$ws.Range("D43").Value = "'0.007741"
$ws.Range("E43").Value = "'1.15%"
$ws.Range("D44").Value = "'0.002238"
$ws.Range("E44").Value = "'-2.52%"
$ws.Range("E45").Value = "'9.40%"
$ws.Range("D46").Value = "'0.00006319"
$ws.Range("E46").Value = "'4.28%"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("D48").Value = "'0.0005808"
$ws.Range("E48").Value = "'0.12%"
$ws.Range("D49").Value = "'24.34"
$ws.Range("E49").Value = "'180.66%"
$ws.Range("D50").Value = "'0.002003"
$ws.Range("E50").Value = "'-25.79%"
$ws.Range("E51").Value = "'0.21%"
